$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Update the statistics (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes) for the countries whose
# figures changed in this data refresh. Columns B..H map to those 7 metrics.
# ---------------------------------------------------------------------------

$updates = @{
    4   = @(1347411, 102,  238080, 1029287, 16816,   7, 80044)   # Estados Unidos
    13  = @(107603,  1383, 86143,  14820,   2675,   51, 6640)    # Iran
    19  = @(42627,   245,  0,      36937,   541,    18, 5440)    # Paises Bajos
    61  = @(4856,    82,   2065,   2783,    2,       0, 8)       # Barein
    64  = @(4402,    369,  558,    3724,    7,       5, 120)     # Afganistan
    77  = @(2117,    27,   1106,   904,     4,       5, 107)     # Bosnia y Herzegovina
    79  = @(1955,    34,   444,    1421,    56,      0, 90)      # Bulgaria
    103 = @(868,     12,   650,    187,     7,       0, 31)      # Albania
    106 = @(845,     36,   234,    585,     4,       0, 26)      # Libano
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    for ($col = 0; $col -lt $vals.Length; $col++) {
        $ws.Cells.Item($row, 2 + $col).Value = $vals[$col]
    }
}

# ---------------------------------------------------------------------------
# The table (rows 4..218) is kept sorted by "Casos totales" (column B) in
# descending order. Re-apply that sort now that several totals changed, so
# the affected countries move to their new rank.
# ---------------------------------------------------------------------------

$dataRange = $ws.Range("A4:H218")
$sortKey = $ws.Range("B4:B218")
$dataRange.Sort($sortKey, 2)

# ---------------------------------------------------------------------------
# "Butan" and "Islas Virgenes Britanicas" are tied on Casos totales (7) and
# swap their relative order in the source refresh even though neither of
# their values changed. Locate that pair (it lands at rows 212/213) and flip
# it explicitly to match the published update.
# ---------------------------------------------------------------------------

$rowButan = $null
$rowVirgenes = $null
for ($r = 4; $r -le 218; $r++) {
    $name = $ws.Cells.Item($r, 1).Value2
    if ($name -eq "Butan") { $rowButan = $r }
    if ($name -eq "Islas Virgenes Britanicas") { $rowVirgenes = $r }
}

# The target order is: Islas Virgenes Britanicas immediately followed by
# Butan. Swap the two rows whenever they are adjacent but not already in
# that order (covers either adjacent ordering produced by the sort).
$needSwap = $false
$topRow = $null
$bottomRow = $null
if ($rowButan -ne $null -and $rowVirgenes -ne $null -and ($rowVirgenes -eq $rowButan + 1)) {
    $needSwap = $true
    $topRow = $rowButan
    $bottomRow = $rowVirgenes
} elseif ($rowButan -ne $null -and $rowVirgenes -ne $null -and ($rowButan -eq $rowVirgenes + 1)) {
    $needSwap = $false
} 

if ($needSwap) {
    $topVals = @()
    $bottomVals = @()
    for ($c = 1; $c -le 8; $c++) {
        $topVals += $ws.Cells.Item($topRow, $c).Value2
        $bottomVals += $ws.Cells.Item($bottomRow, $c).Value2
    }
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($topRow, $c).Value = $bottomVals[$c - 1]
        $ws.Cells.Item($bottomRow, $c).Value = $topVals[$c - 1]
    }
}
